$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.474.37"
$ws.Range("E2").Value = "  -0.26%  "

# Row 3
$ws.Range("D3").Value = "1.838.81"
$ws.Range("E3").Value = "  -0.54%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.33"
$ws.Range("E5").Value = "  -0.80%  "

# Row 6
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5362"
$ws.Range("E7").Value = "  +2.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2934"
$ws.Range("E8").Value = "  -9.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06971"
$ws.Range("E9").Value = "  +2.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.29"
$ws.Range("E10").Value = "  -8.70%  "

# Row 11
$ws.Range("D11").Value = "1.849.57"
$ws.Range("E11").Value = "  -1.80%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7240"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07188"
$ws.Range("E13").Value = "  -7.36%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.10"
$ws.Range("E14").Value = "  +0.60%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.04%  "

# Row 17
$ws.Range("E17").Value = "  -1.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007885"
$ws.Range("E19").Value = "  -0.86%  "

# Row 20
$ws.Range("D20").Value = "26.482.14"
$ws.Range("E20").Value = "  -0.36%  "

# Row 21
$ws.Range("D21").Value = "2.080.22"
$ws.Range("E21").Value = "  -0.92%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.582"
$ws.Range("E22").Value = "  -1.29%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.978"
$ws.Range("E23").Value = "  -0.33%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.163"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.77"
$ws.Range("E25").Value = "  -0.21%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.153"
$ws.Range("E26").Value = "  -0.59%  "

# Row 27
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.708"
$ws.Range("E27").Value = "  +1.68%  "

# Row 28
$ws.Range("E28").Value = "  -0.52%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.63"
$ws.Range("E29").Value = "  -1.31%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.257"
$ws.Range("E30").Value = "  +1.86%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08883"
$ws.Range("E31").Value = "  +1.79%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.029"
$ws.Range("E32").Value = "  -1.72%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04835"
$ws.Range("E33").Value = "  -0.68%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.905"
$ws.Range("E34").Value = "  +1.06%  "

# Row 35
$ws.Range("E35").Value = "  +0.70%  "

# Row 36
$ws.Range("E36").Value = "  -0.14%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.097"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.284"
$ws.Range("E38").Value = "  +0.30%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01705"
$ws.Range("E39").Value = "  -4.41%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4664"
$ws.Range("E40").Value = "  -3.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8993"
$ws.Range("E41").Value = "  -0.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.20"
$ws.Range("E42").Value = "  -3.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.867"
$ws.Range("E43").Value = "  -1.48%  "

# Row 44
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("E45").Value = "  -3.79%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.016"
$ws.Range("E46").Value = "  +0.51%  "

# Row 47
$ws.Range("E47").Value = "  +0.66%  "

# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4044"
$ws.Range("E48").Value = "  -2.92%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.74"
$ws.Range("E49").Value = "  -1.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8921"
$ws.Range("E50").Value = "  -0.52%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05740"
$ws.Range("E51").Value = "  -2.26%  "
